# Daily automated price-data update: prepend a new "today" row (pushing the
# history down) and keep the trailing history row that rolls into view at
# the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (row 1 is the header), shifting all the
# existing date rows down by one.
$ws.Rows.Item(2).Insert()

# New top row: today's date with the same (unchanged) price figures.
# The leading apostrophe forces the date-looking text to stay text instead
# of being auto-converted to a date serial number; ClearFormats() drops the
# quote-prefix style afterwards so the cell matches the plain formatting
# used by the rest of the column.
$ws.Cells.Item(2, 1).Value = "'2026-02-16"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610

# The row that now appears at the bottom (row 89) continues the series with
# the next older date, carrying the same price figures.
$ws.Cells.Item(89, 1).Value = "'2025-11-21"
$ws.Cells.Item(89, 1).ClearFormats()
$ws.Cells.Item(89, 2).Value = 783.5
$ws.Cells.Item(89, 3).Value = 1112
$ws.Cells.Item(89, 4).Value = 3610
